$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 925.2143
$ws.Range("J2").Value = 2042.5
$ws.Range("L2").Value = 2042.5
$ws.Range("N2").Value = -2268.5

$ws.Range("H33").Value = 881.7778
$ws.Range("I33").Value = 854.7857
$ws.Range("J33").Value = 976.25
$ws.Range("K33").Value = 854.7857
$ws.Range("L33").Value = 976.25
$ws.Range("M33").Value = -625.7857
$ws.Range("N33").Value = -1434.25

$ws.Range("H76").Value = 45460548
$ws.Range("I76").Value = 66672030
$ws.Range("J76").Value = 7363.857
$ws.Range("K76").Value = 66672030
$ws.Range("L76").Value = 7363.857
$ws.Range("M76").Value = -66671715
$ws.Range("N76").Value = -7993.857

$ws.Range("H79").Value = 45460548
$ws.Range("I79").Value = 66672030
$ws.Range("J79").Value = 7363.857
$ws.Range("K79").Value = 66672030
$ws.Range("L79").Value = 7363.857
$ws.Range("M79").Value = -66670938
$ws.Range("N79").Value = -9547.857

$ws.Range("H82").Value = 3552.5
$ws.Range("I82").Value = 3552.5
$ws.Range("K82").Value = 10657.5
$ws.Range("M82").Value = -10251.5

$ws.Range("H85").Value = 3552.5
$ws.Range("I85").Value = 3552.5
$ws.Range("K85").Value = 10657.5
$ws.Range("M85").Value = -9253.5

$ws.Range("H86").Value = 41668670
$ws.Range("I86").Value = 83335730
$ws.Range("J86").Value = 1611.3334
$ws.Range("K86").Value = 83335730
$ws.Range("L86").Value = 1611.3334
$ws.Range("M86").Value = -83334607
$ws.Range("N86").Value = -3857.3334

$ws.Range("H89").Value = 41668670
$ws.Range("I89").Value = 83335730
$ws.Range("J89").Value = 1611.3334
$ws.Range("K89").Value = 416678650
$ws.Range("L89").Value = 8056.666999999999
$ws.Range("M89").Value = -416673034
$ws.Range("N89").Value = -19288.667

$ws.Range("H106").Value = 4376.385
$ws.Range("J106").Value = 3879.8
$ws.Range("L106").Value = 3879.8
$ws.Range("N106").Value = -5141.8

$ws.Range("H129").Value = 37053620
$ws.Range("I129").Value = 5953
$ws.Range("J129").Value = 55577456
$ws.Range("K129").Value = 17859
$ws.Range("L129").Value = 166732368
$ws.Range("M129").Value = -12859
$ws.Range("N129").Value = -166742368

$ws.Range("H138").Value = 5216.8696
$ws.Range("I138").Value = 874.2
$ws.Range("J138").Value = 6423.1665
$ws.Range("K138").Value = 2622.6
$ws.Range("L138").Value = 19269.4995
$ws.Range("M138").Value = 2517.4
$ws.Range("N138").Value = -29549.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3006715
$ws.Range("I32").Value = 1705.0555
$ws.Range("J32").Value = 5853566.5
$ws.Range("K32").Value = 1705.0555
$ws.Range("L32").Value = 5853566.5
$ws.Range("M32").Value = -1418.0555
$ws.Range("N32").Value = -5854140.5

$ws.Range("H132").Value = 3907753.5
$ws.Range("I132").Value = 15153640
$ws.Range("J132").Value = 159125
$ws.Range("K132").Value = 45460920
$ws.Range("L132").Value = 477375
$ws.Range("M132").Value = -45458390
$ws.Range("N132").Value = -482435

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4762.4814
$ws.Range("I99").Value = 2347.1052
$ws.Range("K99").Value = 2347.1052
$ws.Range("M99").Value = -849.1052

$ws.Range("H134").Value = 2515733.5
$ws.Range("I134").Value = 2811184.5
$ws.Range("J134").Value = 4400
$ws.Range("K134").Value = 8433553.5
$ws.Range("L134").Value = 13200
$ws.Range("M134").Value = -8431018.5
$ws.Range("N134").Value = -18270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6426.8335
$ws.Range("I31").Value = 1832.4286
$ws.Range("J31").Value = 12859
$ws.Range("K31").Value = 1832.4286
$ws.Range("L31").Value = 12859
$ws.Range("M31").Value = -1537.4286
$ws.Range("N31").Value = -13449

$ws.Range("H34").Value = 6426.8335
$ws.Range("I34").Value = 1832.4286
$ws.Range("J34").Value = 12859
$ws.Range("K34").Value = 1832.4286
$ws.Range("L34").Value = 12859
$ws.Range("M34").Value = -1630.4286
$ws.Range("N34").Value = -13263

$ws.Range("H58").Value = 35721160
$ws.Range("I58").Value = 52636320
$ws.Range("J58").Value = 11378.223
$ws.Range("K58").Value = 52636320
$ws.Range("L58").Value = 11378.223
$ws.Range("M58").Value = -52636117
$ws.Range("N58").Value = -11784.223

$ws.Range("H62").Value = 8132.2
$ws.Range("I62").Value = 8132.2
$ws.Range("K62").Value = 8132.2
$ws.Range("M62").Value = -7508.2

$ws.Range("H65").Value = 8132.2
$ws.Range("I65").Value = 8132.2
$ws.Range("K65").Value = 40661
$ws.Range("M65").Value = -37541

$ws.Range("H132").Value = 5912.1577
$ws.Range("I132").Value = 5133.3335
$ws.Range("K132").Value = 15400.0005
$ws.Range("M132").Value = -12870.0005

$ws.Range("H136").Value = 35721160
$ws.Range("I136").Value = 52636320
$ws.Range("J136").Value = 11378.223
$ws.Range("K136").Value = 157908960
$ws.Range("L136").Value = 34134.669
$ws.Range("M136").Value = -157906410
$ws.Range("N136").Value = -39234.669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 44945.855
$ws.Range("J121").Value = 44945.855
$ws.Range("L121").Value = 134837.565
$ws.Range("N121").Value = -137457.565

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7172.364
$ws.Range("I102").Value = 6458.9473
$ws.Range("J102").Value = 8140.5713
$ws.Range("K102").Value = 6458.9473
$ws.Range("L102").Value = 8140.5713
$ws.Range("M102").Value = -4836.9473
$ws.Range("N102").Value = -11384.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1954.3235
$ws.Range("I55").Value = 704.46155
$ws.Range("J55").Value = 2728.0476
$ws.Range("K55").Value = 704.46155
$ws.Range("L55").Value = 2728.0476
$ws.Range("M55").Value = -531.46155
$ws.Range("N55").Value = -3074.0476

$ws.Range("H68").Value = 2124.6155
$ws.Range("I68").Value = 2138.182
$ws.Range("J68").Value = 2050
$ws.Range("K68").Value = 2138.182
$ws.Range("L68").Value = 2050
$ws.Range("M68").Value = -1389.182
$ws.Range("N68").Value = -3548

$ws.Range("H71").Value = 2124.6155
$ws.Range("I71").Value = 2138.182
$ws.Range("J71").Value = 2050
$ws.Range("K71").Value = 10690.91
$ws.Range("L71").Value = 10250
$ws.Range("M71").Value = -6946.91
$ws.Range("N71").Value = -17738

$ws.Range("H82").Value = 2769.2068
$ws.Range("I82").Value = 1553.1
$ws.Range("J82").Value = 5471.6665
$ws.Range("K82").Value = 1553.1
$ws.Range("L82").Value = 5471.6665
$ws.Range("M82").Value = -1192.1
$ws.Range("N82").Value = -6193.6665

$ws.Range("H85").Value = 2769.2068
$ws.Range("I85").Value = 1553.1
$ws.Range("J85").Value = 5471.6665
$ws.Range("K85").Value = 1553.1
$ws.Range("L85").Value = 5471.6665
$ws.Range("M85").Value = -305.0999999999999
$ws.Range("N85").Value = -7967.6665

$ws.Range("H132").Value = 1729.9714
$ws.Range("I132").Value = 1616.6364
$ws.Range("K132").Value = 4849.9092
$ws.Range("M132").Value = -2319.9092

$ws.Range("H136").Value = 9808587
$ws.Range("I136").Value = 35720116
$ws.Range("K136").Value = 107160348
$ws.Range("M136").Value = -107157798

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17016.572
$ws.Range("I62").Value = 11185.625
$ws.Range("J62").Value = 24791.166
$ws.Range("K62").Value = 11185.625
$ws.Range("L62").Value = 24791.166
$ws.Range("M62").Value = -10561.625
$ws.Range("N62").Value = -26039.166

$ws.Range("H65").Value = 17016.572
$ws.Range("I65").Value = 11185.625
$ws.Range("J65").Value = 24791.166
$ws.Range("K65").Value = 55928.125
$ws.Range("L65").Value = 123955.83
$ws.Range("M65").Value = -52808.125
$ws.Range("N65").Value = -130195.83

$ws.Range("H107").Value = 456.0625
$ws.Range("I107").Value = 584.9
$ws.Range("J107").Value = 241.33333
$ws.Range("K107").Value = 1754.7
$ws.Range("L107").Value = 723.99999
$ws.Range("M107").Value = 165.3000000000002
$ws.Range("N107").Value = -4563.99999

$ws.Range("H132").Value = 5903.684
$ws.Range("I132").Value = 3618.0688
$ws.Range("K132").Value = 10854.2064
$ws.Range("M132").Value = -8324.206399999999
